# Auto-generated: applies cryptos list price/volume update (2023-12-20 GitHub Actions run)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "42.848.68"
$ws.Range("D3").Value = "2.210.65"
$ws.Range("E3").Value = "  -1.62%  "
$ws.Range("E4").Value = "  -0.11%  "
$ws.Range("D5").Value = "'257.11"
$ws.Range("E5").Value = "  +5.02%  "
$ws.Range("D6").Value = "'0.616"
$ws.Range("E6").Value = "  -0.14%  "
$ws.Range("D7").Value = "'76.88"
$ws.Range("E7").Value = "  +1.79%  "
$ws.Range("E8").Value = "  +0.00%  "
$ws.Range("D9").Value = "'0.594"
$ws.Range("E9").Value = "  -2.30%  "
$ws.Range("D10").Value = "'42.16"
$ws.Range("E10").Value = "  +2.90%  "
$ws.Range("D11").Value = "'0.0907"
$ws.Range("E11").Value = "  -2.83%  "
$ws.Range("E12").Value = "  -0.43%  "
$ws.Range("E13").Value = "  +1.01%  "
$ws.Range("D14").Value = "2.543.18"
$ws.Range("E14").Value = "  -1.60%  "
$ws.Range("D15").Value = "'14.48"
$ws.Range("E15").Value = "  -1.01%  "
$ws.Range("D16").Value = "2.221.25"
$ws.Range("E16").Value = "  -1.03%  "
$ws.Range("D17").Value = "'0.783"
$ws.Range("E17").Value = "  -1.83%  "
$ws.Range("D18").Value = "42.830.42"
$ws.Range("E18").Value = "  -0.43%  "
$ws.Range("E19").Value = "  -2.20%  "
$ws.Range("D20").Value = "'71.15"
$ws.Range("E20").Value = "  -0.14%  "
$ws.Range("D21").Value = "'5.97"
$ws.Range("E21").Value = "  -0.21%  "
$ws.Range("D22").Value = "'2.21"
$ws.Range("E22").Value = "  +0.34%  "
$ws.Range("D23").Value = "'230.22"
$ws.Range("E23").Value = "  -0.06%  "
$ws.Range("D24").Value = "'9.37"
$ws.Range("E24").Value = "  -6.70%  "
$ws.Range("E25").Value = "  -0.12%  "
$ws.Range("D26").Value = "'43.11"
$ws.Range("E26").Value = "  +11.28%  "
$ws.Range("D27").Value = "'10.74"
$ws.Range("E27").Value = "  -1.31%  "
$ws.Range("D28").Value = "'3.34"
$ws.Range("E28").Value = "  -3.45%  "
$ws.Range("E29").Value = "  -1.76%  "
$ws.Range("E30").Value = "  +0.09%  "
$ws.Range("D31").Value = "'173.08"
$ws.Range("E31").Value = "  -0.19%  "
$ws.Range("D32").Value = "'20.33"
$ws.Range("E32").Value = "  +0.07%  "
$ws.Range("D33").Value = "'0.0872"
$ws.Range("E33").Value = "  +9.44%  "
$ws.Range("E34").Value = "  -1.77%  "
$ws.Range("E35").Value = "  -0.12%  "
$ws.Range("D36").Value = "'0.0363"
$ws.Range("E36").Value = "  +8.02%  "
$ws.Range("E37").Value = "  -2.39%  "
$ws.Range("D38").Value = "'4.35"
$ws.Range("E38").Value = "  -0.13%  "
$ws.Range("D39").Value = "'12.87"
$ws.Range("E39").Value = "  -1.78%  "
$ws.Range("B40").Value = "LidoDAOToken"
$ws.Range("C40").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D40").Value = "'2.11"
$ws.Range("E40").Value = "  -1.16%  "
$ws.Range("B41").Value = "NEARProtocol"
$ws.Range("C41").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D41").Value = "'2.80"
$ws.Range("E41").Value = "  +17.41%  "
$ws.Range("E42").Value = "  -2.49%  "
$ws.Range("E43").Value = "  -3.83%  "
$ws.Range("D44").Value = "'60.00"
$ws.Range("E44").Value = "  +0.38%  "
$ws.Range("D45").Value = "'102.71"
$ws.Range("E45").Value = "  -2.61%  "
$ws.Range("E46").Value = "  -4.32%  "
$ws.Range("D47").Value = "'0.0978"
$ws.Range("E47").Value = "  -1.56%  "
$ws.Range("D48").Value = "'0.461"
$ws.Range("E48").Value = "  -5.49%  "
$ws.Range("E49").Value = "  +0.49%  "
$ws.Range("E50").Value = "  -0.93%  "
$ws.Range("B51").Value = "RocketPoolETH"
$ws.Range("C51").Value = "https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth"
$ws.Range("D51").Value = "2.431.12"
$ws.Range("E51").Value = "  -1.15%  "
